$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues / xlPasteFormats constants (as used by Excel COM PasteSpecial)
$xlPasteValues = -4163
$xlPasteFormats = -4122

# 1) Shift the old row 3 (claim 2 / preproduccion / tcorvetto) down to row 5,
#    preserving both its values and its cell formatting (quote-prefixed text style).
$ws.Range("A3:F3").Copy()
$ws.Range("A5:F5").PasteSpecial($xlPasteValues)
$ws.Range("A3:F3").Copy()
$ws.Range("A5:F5").PasteSpecial($xlPasteFormats)

# 2) Clear out the old row 3 content now that it lives on in row 5.
$ws.Range("A3:F3").ClearContents()

# 3) Build new row 3 - same environment/user as row 2 - by copying row 2's
#    B:E cells (values + formats) straight across.
$ws.Range("B2:E2").Copy()
$ws.Range("B3:E3").PasteSpecial($xlPasteValues)
$ws.Range("B2:E2").Copy()
$ws.Range("B3:E3").PasteSpecial($xlPasteFormats)

# 4) Build new row 4 the same way.
$ws.Range("B2:E2").Copy()
$ws.Range("B4:E4").PasteSpecial($xlPasteValues)
$ws.Range("B2:E2").Copy()
$ws.Range("B4:E4").PasteSpecial($xlPasteFormats)

# 5) Give F3/F4 the same quote-prefixed text formatting as F5 (carried over
#    from the original row 3), then fill in the three claim-number values in
#    the exact order the new shared strings must be introduced: F3's number
#    first, then F2's (new) number, then F4's number.
$ws.Range("F5").Copy()
$ws.Range("F3").PasteSpecial($xlPasteFormats)
$ws.Range("F5").Copy()
$ws.Range("F4").PasteSpecial($xlPasteFormats)

$ws.Range("F3").Value = "'1120194100412"
$ws.Range("F2").Value = "1220194200667"
$ws.Range("F4").Value = "'0420194406717"

$ws.Range("F5").Select()
